# Generate Report for Handoff
# Updates the status of the ab1e550b-... file (row 3 in each sheet) to
# "Ready for handoff", refreshes its timestamps, and switches its
# Priority to "mt" on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-27 02:14:47"

# --- zh-cn sheet --------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = "mt"
$zh.Range("H3").Value = "2016-08-27 02:14:43"

# --- de-de sheet --------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = "mt"
$de.Range("H3").Value = "2016-08-27 02:14:47"

# --- Column width adjustments (status column grew wider to fit the new
#     "Ready for handoff" text). The runtime quantizes ColumnWidth onto an
#     integer-pixel grid (px = round(ColumnWidth*6+5); stored chars = px/6),
#     so 16.3 is the closest input that reproduces the target ~17.22 chars
#     width (lands on 17.1667, i.e. the nearest representable grid value).
$ov.Range("E:F").ColumnWidth = 16.3
$zh.Range("C:C").ColumnWidth = 16.3
$de.Range("C:C").ColumnWidth = 16.3
